# Adds the missing daily metrics (columns B, C, E, F, G, H, I, J) for rows
# 54-62 (dates 2025-12-01 .. 2025-12-09, i.e. serials 46011-46019). Column D
# already holds a "shared" formula (B-C) inherited from the D35:D63 fill, so
# it recalculates automatically once B/C are populated - no need to touch it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 54 (2025-12-01)
$ws.Range("B54").Value = 1213
$ws.Range("C54").Value = 843
$ws.Range("E54").Value = 18
$ws.Range("F54").Value = 1
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 138
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0

# Row 55 (2025-12-02)
$ws.Range("B55").Value = 1294
$ws.Range("C55").Value = 405
$ws.Range("E55").Value = 14
$ws.Range("F55").Value = 1
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 108
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0

# Row 56 (2025-12-03)
$ws.Range("B56").Value = 2735
$ws.Range("C56").Value = 2204
$ws.Range("E56").Value = 55
$ws.Range("F56").Value = 5
$ws.Range("G56").Value = 5
$ws.Range("H56").Value = 254
$ws.Range("I56").Value = 6
$ws.Range("J56").Value = 36

# Row 57 (2025-12-04)
$ws.Range("B57").Value = 2573
$ws.Range("C57").Value = 1882
$ws.Range("E57").Value = 60
$ws.Range("F57").Value = 5
$ws.Range("G57").Value = 5
$ws.Range("H57").Value = 253
$ws.Range("I57").Value = 10
$ws.Range("J57").Value = 178

# Row 58 (2025-12-05)
$ws.Range("B58").Value = 3147
$ws.Range("C58").Value = 2706
$ws.Range("E58").Value = 40
$ws.Range("F58").Value = 1
$ws.Range("G58").Value = 2
$ws.Range("H58").Value = 314
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0

# Row 59 (2025-12-06)
$ws.Range("B59").Value = 2227
$ws.Range("C59").Value = 1735
$ws.Range("E59").Value = 46
$ws.Range("F59").Value = 3
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = 251
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0

# Row 60 (2025-12-07)
$ws.Range("B60").Value = 2596
$ws.Range("C60").Value = 2247
$ws.Range("E60").Value = 62
$ws.Range("F60").Value = 6
$ws.Range("G60").Value = 6
$ws.Range("H60").Value = 259
$ws.Range("I60").Value = 26
$ws.Range("J60").Value = 18

# Row 61 (2025-12-08)
$ws.Range("B61").Value = 2720
$ws.Range("C61").Value = 2370
$ws.Range("E61").Value = 47
$ws.Range("F61").Value = 6
$ws.Range("G61").Value = 7
$ws.Range("H61").Value = 142
$ws.Range("I61").Value = 6
$ws.Range("J61").Value = 0

# Row 62 (2025-12-09)
$ws.Range("B62").Value = 591
$ws.Range("C62").Value = 288
$ws.Range("E62").Value = 21
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 107
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0

# Match the author's final cursor position/selection (G54) from the diff.
$ws.Range("G54").Select()
